$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2024-09-29 Sunday" "2024-09-30 Monday"

Replace-Text "53÷6=8, 5" "71÷9=7, 8"
Replace-Text "78÷7=11, 1" "16÷6=2, 4"
Replace-Text "63÷4=15, 3" "98÷6=16, 2"
Replace-Text "77÷9=8, 5" "65÷8=8, 1"
Replace-Text "25÷6=4, 1" "13÷5=2, 3"

Replace-Text "54÷4=13, 2" "41÷8=5, 1"
Replace-Text "32÷6=5, 2" "64÷5=12, 4"
Replace-Text "49÷8=6, 1" "93÷5=18, 3"
Replace-Text "76÷2=38, 0" "95÷6=15, 5"
Replace-Text "47÷3=15, 2" "59÷5=11, 4"

Replace-Text "20÷3=6, 2" "16÷5=3, 1"
Replace-Text "13÷8=1, 5" "74÷5=14, 4"
Replace-Text "26÷8=3, 2" "50÷8=6, 2"
Replace-Text "63÷8=7, 7" "65÷6=10, 5"
Replace-Text "86÷9=9, 5" "63÷5=12, 3"

Replace-Text "35÷4=8, 3" "63÷6=10, 3"
Replace-Text "28÷5=5, 3" "15÷5=3, 0"
Replace-Text "77÷5=15, 2" "44÷3=14, 2"
Replace-Text "11÷3=3, 2" "63÷9=7, 0"
Replace-Text "10÷5=2, 0" "31÷8=3, 7"

Replace-Text "28÷3=9, 1" "89÷8=11, 1"
Replace-Text "66÷7=9, 3" "63÷7=9, 0"
Replace-Text "88÷3=29, 1" "19÷4=4, 3"
Replace-Text "69÷2=34, 1" "47÷5=9, 2"
Replace-Text "99÷5=19, 4" "94÷4=23, 2"

Write-Output "Done"
